$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching style of existing header cells (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-93
$data = @(
    @{Row=2; I=8; J=8},
    @{Row=3; I=5; J=5},
    @{Row=4; I=5; J=7},
    @{Row=5; I=8; J=8},
    @{Row=6; I=8; J=8},
    @{Row=7; I=8; J=9},
    @{Row=8; I=8; J=8},
    @{Row=9; I=7; J=7},
    @{Row=10; I=8; J=8},
    @{Row=11; I=9; J=9},
    @{Row=12; I=9; J=9},
    @{Row=13; I=8; J=8},
    @{Row=14; I=9; J=9},
    @{Row=15; I=8; J=8},
    @{Row=16; I=7; J=7},
    @{Row=17; I=7; J=7},
    @{Row=18; I=8; J=8},
    @{Row=19; I=8; J=8},
    @{Row=20; I=8; J=8},
    @{Row=21; I=8; J=8},
    @{Row=22; I=8; J=8},
    @{Row=23; I=8; J=8},
    @{Row=24; I=9; J=9},
    @{Row=25; I=8; J=8},
    @{Row=26; I=8; J=8},
    @{Row=27; I=8; J=8},
    @{Row=28; I=8; J=8},
    @{Row=29; I=8; J=8},
    @{Row=30; I=8; J=8},
    @{Row=31; I=8; J=8},
    @{Row=32; I=8; J=8},
    @{Row=33; I=9; J=9},
    @{Row=34; I=8; J=8},
    @{Row=35; I=9; J=9},
    @{Row=36; I=8; J=8},
    @{Row=37; I=8; J=8},
    @{Row=38; I=8; J=8},
    @{Row=39; I=8; J=8},
    @{Row=40; I=8; J=8},
    @{Row=41; I=9; J=9},
    @{Row=42; I=9; J=9},
    @{Row=43; I=9; J=9},
    @{Row=44; I=8; J=8},
    @{Row=45; I=8; J=8},
    @{Row=46; I=8; J=8},
    @{Row=47; I=9; J=9},
    @{Row=48; I=8; J=8},
    @{Row=49; I=9; J=9},
    @{Row=50; I=8; J=8},
    @{Row=51; I=8; J=8},
    @{Row=52; I=8; J=8},
    @{Row=53; I=9; J=9},
    @{Row=54; I=8; J=8},
    @{Row=55; I=6; J=6},
    @{Row=56; I=9; J=10},
    @{Row=57; I=9; J=9},
    @{Row=58; I=9; J=9},
    @{Row=59; I=9; J=9},
    @{Row=60; I=9; J=9},
    @{Row=61; I=9; J=9},
    @{Row=62; I=9; J=9},
    @{Row=63; I=9; J=9},
    @{Row=64; I=10; J=10},
    @{Row=65; I=9; J=9},
    @{Row=66; I=8; J=8},
    @{Row=67; I=8; J=8},
    @{Row=68; I=8; J=8},
    @{Row=69; I=6; J=6},
    @{Row=70; I=8; J=8},
    @{Row=71; I=8; J=8},
    @{Row=72; I=9; J=9},
    @{Row=73; I=8; J=8},
    @{Row=74; I=8; J=8},
    @{Row=75; I=8; J=8},
    @{Row=76; I=9; J=9},
    @{Row=77; I=7; J=8},
    @{Row=78; I=8; J=8},
    @{Row=79; I=8; J=8},
    @{Row=80; I=7; J=8},
    @{Row=81; I=8; J=8},
    @{Row=82; I=8; J=8},
    @{Row=83; I=8; J=8},
    @{Row=84; I=8; J=8},
    @{Row=85; I=8; J=8},
    @{Row=86; I=7; J=7},
    @{Row=87; I=7; J=7},
    @{Row=88; I=5; J=5},
    @{Row=89; I=7; J=7},
    @{Row=90; I=6; J=6},
    @{Row=91; I=3; J=3},
    @{Row=92; I=3; J=3},
    @{Row=93; I=3; J=3}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 9).Value = $row.I   # column I
    $ws.Cells.Item($row.Row, 10).Value = $row.J  # column J
}
